# Update attendance summary cells from 0 to 1 for each date row (rows 3-18).
# Each row gets its own subset of columns (D, E, G, H) flipped to 1, matching
# the specific attendance classification (Real/Duplicate/Invalid/Absent)
# recorded for that date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("D", "E")
    6  = @("D", "E")
    7  = @("H")
    8  = @("H")
    9  = @("H")
    10 = @("D", "E")
    11 = @("H")
    12 = @("H")
    13 = @("D", "E")
    14 = @("G", "H")
    15 = @("G", "H")
    16 = @("H")
    17 = @("D", "E")
    18 = @("G", "H")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
